# finish rate algorithm demo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CONDITION cell for the worker_connections rule from an
# equality check to a less-than-or-equal check.
$ws.Range("C8").Value = 'worker_connections<=$param'

# Update the saved selection to C14 (matches the author's final cursor
# position before saving).
$ws.Range("C14").Select()
